$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E9: change from text "*" to numeric 0.5
$ws.Range("E9").Value = 0.5

# E13: change from text "*socket.io" to text "*"
$ws.Range("E13").Value = "*"

# E15: add new numeric cell
$ws.Range("E15").Value = 0.5

# Update view: topLeftCell A4 -> A6, selection B13 -> B16
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B16").Select()
